$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1642.8
$ws.Range("J121").Value = 1642.8
$ws.Range("L121").Value = 4928.4
$ws.Range("N121").Value = -8422.4
$ws.Range("H132").Value = 1279.8462
$ws.Range("I132").Value = 1169.8334
$ws.Range("K132").Value = 3509.5002
$ws.Range("M132").Value = -979.5001999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1404
$ws.Range("I2").Value = 1415.2
$ws.Range("J2").Value = 1348
$ws.Range("K2").Value = 1415.2
$ws.Range("L2").Value = 1348
$ws.Range("M2").Value = -1302.2
$ws.Range("N2").Value = -1574
$ws.Range("H14").Value = 1006
$ws.Range("I14").Value = 1006
$ws.Range("K14").Value = 1006
$ws.Range("M14").Value = -831
$ws.Range("H61").Value = 4994.8184
$ws.Range("I61").Value = 5368.5
$ws.Range("J61").Value = 3998.3333
$ws.Range("K61").Value = 5368.5
$ws.Range("L61").Value = 3998.3333
$ws.Range("M61").Value = -5156.5
$ws.Range("N61").Value = -4422.3333
$ws.Range("H116").Value = 1404
$ws.Range("I116").Value = 1415.2
$ws.Range("J116").Value = 1348
$ws.Range("K116").Value = 1415.2
$ws.Range("L116").Value = 1348
$ws.Range("M116").Value = 878.8
$ws.Range("N116").Value = -5936
$ws.Range("H122").Value = 2830.1667
$ws.Range("I122").Value = 2797.5
$ws.Range("K122").Value = 8392.5
$ws.Range("M122").Value = -5942.5
$ws.Range("H136").Value = 4994.8184
$ws.Range("I136").Value = 5368.5
$ws.Range("J136").Value = 3998.3333
$ws.Range("K136").Value = 16105.5
$ws.Range("L136").Value = 11994.9999
$ws.Range("M136").Value = -13555.5
$ws.Range("N136").Value = -17094.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1404
$ws.Range("I3").Value = 1415.2
$ws.Range("J3").Value = 1348
$ws.Range("K3").Value = 1415.2
$ws.Range("L3").Value = 1348
$ws.Range("M3").Value = -1301.2
$ws.Range("N3").Value = -1576
$ws.Range("H134").Value = 2079.4211
$ws.Range("I134").Value = 1824
$ws.Range("K134").Value = 5472
$ws.Range("M134").Value = -2937

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6459
$ws.Range("J31").Value = 9981.666999999999
$ws.Range("L31").Value = 9981.666999999999
$ws.Range("N31").Value = -10571.667
$ws.Range("H34").Value = 6459
$ws.Range("J34").Value = 9981.666999999999
$ws.Range("L34").Value = 9981.666999999999
$ws.Range("N34").Value = -10385.667
$ws.Range("H94").Value = 2272.25
$ws.Range("J94").Value = 2196.6667
$ws.Range("L94").Value = 2196.6667
$ws.Range("N94").Value = -3098.6667
$ws.Range("H99").Value = 3669.4285
$ws.Range("I99").Value = 3624.25
$ws.Range("K99").Value = 3624.25
$ws.Range("M99").Value = -2126.25
$ws.Range("H126").Value = 3669.4285
$ws.Range("I126").Value = 3624.25
$ws.Range("K126").Value = 10872.75
$ws.Range("M126").Value = -8402.75
$ws.Range("H132").Value = 3195.3333
$ws.Range("I132").Value = 2483
$ws.Range("K132").Value = 7449
$ws.Range("M132").Value = -4919

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 96162470
$ws.Range("I4").Value = 61900404
$ws.Range("K4").Value = 185701212
$ws.Range("M4").Value = -185701100
$ws.Range("H5").Value = 5815.1665
$ws.Range("I5").Value = 998
$ws.Range("K5").Value = 2994
$ws.Range("M5").Value = -2882
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -37
$ws.Range("H11").Value = 501
$ws.Range("I11").Value = 501
$ws.Range("K11").Value = 1503
$ws.Range("M11").Value = -1363
$ws.Range("H40").Value = 199.5
$ws.Range("I40").Value = 199.5
$ws.Range("K40").Value = 798
$ws.Range("M40").Value = -729
$ws.Range("H80").Value = 6840.4
$ws.Range("I80").Value = 901
$ws.Range("J80").Value = 10800
$ws.Range("K80").Value = 2703
$ws.Range("L80").Value = 32400
$ws.Range("M80").Value = -1767
$ws.Range("N80").Value = -34272
$ws.Range("H83").Value = 6840.4
$ws.Range("I83").Value = 901
$ws.Range("J83").Value = 10800
$ws.Range("K83").Value = 8109
$ws.Range("L83").Value = 97200
$ws.Range("M83").Value = -3429
$ws.Range("N83").Value = -106560
$ws.Range("H114").Value = 3378.2
$ws.Range("I114").Value = 2187.6667
$ws.Range("J114").Value = 5164
$ws.Range("K114").Value = 6563.000100000001
$ws.Range("L114").Value = 15492
$ws.Range("M114").Value = -3309.000100000001
$ws.Range("N114").Value = -22000
$ws.Range("H122").Value = 3541.0408
$ws.Range("I122").Value = 440.66666
$ws.Range("J122").Value = 3743.239
$ws.Range("K122").Value = 3965.99994
$ws.Range("L122").Value = 33689.151
$ws.Range("M122").Value = -1515.99994
$ws.Range("N122").Value = -38589.151
$ws.Range("H135").Value = 5815.1665
$ws.Range("I135").Value = 998
$ws.Range("K135").Value = 8982
$ws.Range("M135").Value = -6447

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("H46").Value = 2987.5
$ws.Range("I46").Value = 2987.5
$ws.Range("K46").Value = 2987.5
$ws.Range("M46").Value = -2799.5
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13828.833
$ws.Range("I81").Value = 11323
$ws.Range("K81").Value = 22646
$ws.Range("M81").Value = -21585
$ws.Range("H84").Value = 13828.833
$ws.Range("I84").Value = 11323
$ws.Range("K84").Value = 113230
$ws.Range("M84").Value = -107926
$ws.Range("H113").Value = 1199.2
$ws.Range("I113").Value = 1124
$ws.Range("K113").Value = 3372
$ws.Range("M113").Value = -1202
$ws.Range("H126").Value = 2200
$ws.Range("I126").Value = 2200
$ws.Range("K126").Value = 6600
$ws.Range("M126").Value = -4130
$ws.Range("H132").Value = 2964.4348
$ws.Range("I132").Value = 2319.6562
$ws.Range("K132").Value = 6958.9686
$ws.Range("M132").Value = -4428.9686
$ws.Range("H136").Value = 9935.4
$ws.Range("I136").Value = 10673
$ws.Range("J136").Value = 5141
$ws.Range("K136").Value = 32019
$ws.Range("M136").Value = -29469
$ws.Range("N136").Value = -20523
